$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-16
# from serial date 45204 (2023-10-05) to 45207 (2023-10-08).
for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45207
    }
}
